$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the "Example: An MMORPG ..." bullet paragraph under Flyweight.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Example: An MMORPG*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the flyweight 'Example: An MMORPG' paragraph."
}

# ---------------------------------------------------------------------------
# 2) "refer" -> "refer to those names"  (narrow, whole-word Find so only
#    this exact word is touched).
# ---------------------------------------------------------------------------
$findRange = $target.Range.Duplicate
$findRange.Find.Execute("refer", $true, $true, $false, $false, $false, `
                         $true, 0, $false, "refer to those names", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the old " to those names" tail (now redundant) and the
#    trailing "." and replace with ". This concept can be generalised to
#    any data type." The old _GoBack bookmark sits right after "refer";
#    drop it here since it is being relocated into the new paragraph below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$paraText = $target.Range.Text
$tailIdx = $paraText.IndexOf(" to those names.")
$tailStart = $target.Range.Start + $tailIdx
$tailEnd = $target.Range.End - 1
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Text = ". This concept can be generalised to any data type."

# ---------------------------------------------------------------------------
# 4) Insert a brand-new bullet paragraph right after it, matching the
#    CNormal / numPr(ilvl0,numId3) / spacing-after0 formatting.
# ---------------------------------------------------------------------------
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()
$newPara.Range.Text = "To use the pattern: you can use static data members, singleton classes, or a shared object that uniquely stores data of a certain type. Anytime that type is instantiated, you must check to see if the data already exists and store its index/iterator instead of storing the actual data in the client."

# ---------------------------------------------------------------------------
# 5) Re-insert the _GoBack bookmark in the new paragraph, right before
#    " the actual data in the client."
# ---------------------------------------------------------------------------
$newParaText = $newPara.Range.Text
$markerIdx = $newParaText.IndexOf(" the actual data in the client.")
$markerPos = $newPara.Range.Start + $markerIdx
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 6) Drop the stale <w:lastRenderedPageBreak/> cached in the "Google Test"
#    bullet under Extracurricular: round-trip the run text (append then
#    restore) so Word rebuilds the run without the rendering-cache marker,
#    without touching the paragraph mark / paraId.
# ---------------------------------------------------------------------------
$gt = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Google Test can be used*") {
        $gt = $p
        break
    }
}
if ($gt -eq $null) {
    throw "Could not locate the 'Google Test can be used' paragraph."
}
$gtStart = $gt.Range.Start
$gtEnd = $gt.Range.End
$gtRange = $d.Range($gtStart, $gtEnd - 1)
$gtOrig = $gtRange.Text
$gtRange.Text = $gtOrig + "X"
$gtRange2 = $d.Range($gtStart, $gtEnd)
$gtRange2.Text = $gtOrig

Write-Output "Done"
